$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("C1").Value = "sub"
$ws.Range("D1").Value = "instructorId"

$ws.Range("B2").Value = '{"quiz_questions":[{"id":"1","question":"Sparing in words, using very few words","possible_answers":"Loquacious,Laconic,Integral,Judicious","correct_answer":"2"},{"id":"2","question":"Freedom from punishment","possible_answers":"Impunity,Indolent,Jaded,Hedonist","correct_answer":"1"},{"id":"3","question":"An instrument for measuring","possible_answers":"Gallant ,Gauge ,Gamester ,Gastric ","correct_answer":"2"},{"id":"4","question":"Inflammation of the stomach  ","possible_answers":"Insurgent ,Impudence ,Garrulous ,Gastritis ","correct_answer":"4"},{"id”:”5”,”question":"unoriginal","possible_answers":"Imperious ,Impudence ,Hackneyed ,Inane ","correct_answer":"3"},{"id”:”6”,”question":"To draw principle inferences","possible_answers":"Gaseous ,Generalize ,Garrote ,Garrison ","correct_answer":"2"},{"id”:”7”,”question":"To imbue with life or animation ","possible_answers":"Galvanize ,Garnish ,Garrulous ,Gaseous ","correct_answer":"1"},{"id”:”8”,”question":"Decorate or embellish ","possible_answers":"Gamut ,Genealogist ,Gambol ,Garnish ","correct_answer":"4"},{"id”:”9”,”question":"Given to constant trivial talking","possible_answers":"Gaseous ,Gendarme ,Garrulous ,Genitive ","correct_answer":"3"},{"id”:”10”,”question":"To risk money or other possession on an event","possible_answers":"Gallant ,Gaily ,Gamble ,Gastronomy ","correct_answer":"3"}]}'

$ws.Range("C2").Value = "dbms"
$ws.Range("D2").Value = "ShriRadhey"

$null = $ws.Range("A2").Select()
